$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 71 (Oregon Samples tidepool 13, 11/21/2019)
# Copy the date-formatted cell above so the new cell reuses the existing
# date style instead of minting a new number format.
$ws.Range("A70").Copy($ws.Range("A71"))
$ws.Range("A71").Value = 43790

$ws.Range("B71").Value = 2194.45390123274
$ws.Range("C71").Value = 2207.0300000000002
$ws.Range("D71").Formula = "=100*(B71-C71)/C71"
$ws.Range("E71").Value = 169
$ws.Range("F71").Value = "crm opened 11/19/2019"

# Restore the selection left by the editor after entering the row
$ws.Range("C67").Select()
